# Auto-generated edit script applying the Balmung_Profits diff
# Updates cached market-price / profit values across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 22727788

$ws.Range("H83").Value = 22727788

$ws.Range("H111").Value = 1995
$ws.Range("I111").Value = 1995
$ws.Range("K111").Value = 5985
$ws.Range("M111").Value = -2918

$ws.Range("H138").Value = 6618.7812
$ws.Range("J138").Value = 4270.391
$ws.Range("L138").Value = 12811.173
$ws.Range("N138").Value = -23091.173

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 303914.22
$ws.Range("I74").Value = 2093.6309
$ws.Range("K74").Value = 2093.6309
$ws.Range("M74").Value = -1219.6309

$ws.Range("H77").Value = 303914.22
$ws.Range("I77").Value = 2093.6309
$ws.Range("K77").Value = 10468.1545
$ws.Range("M77").Value = -6100.154500000001

$ws.Range("H102").Value = 2644.25
$ws.Range("I102").Value = 2236.2856
$ws.Range("K102").Value = 2236.2856
$ws.Range("M102").Value = -614.2856000000002

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents() | Out-Null

$ws.Range("H132").Value = 3108.25
$ws.Range("I132").Value = 1929.1666
$ws.Range("K132").Value = 5787.4998
$ws.Range("M132").Value = -3257.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7751.516
$ws.Range("I107").Value = 8078.3447
$ws.Range("K107").Value = 8078.3447
$ws.Range("M107").Value = -6158.3447

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 11050
$ws.Range("J21").Value = 21500
$ws.Range("L21").Value = 21500
$ws.Range("N21").Value = -21970

$ws.Range("H134").Value = 2425.2285
$ws.Range("I134").Value = 2295.75
$ws.Range("J134").Value = 2707.7273
$ws.Range("K134").Value = 6887.25
$ws.Range("L134").Value = 8123.1819
$ws.Range("M134").Value = -4352.25
$ws.Range("N134").Value = -13193.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 10499.74
$ws.Range("J76").Value = 10418.102
$ws.Range("L76").Value = 31254.306
$ws.Range("N76").Value = -32020.306

$ws.Range("H79").Value = 10499.74
$ws.Range("J79").Value = 10418.102
$ws.Range("L79").Value = 31254.306
$ws.Range("N79").Value = -33906.306

$ws.Range("H104").Value = 7582.625
$ws.Range("J104").Value = 8094.4287
$ws.Range("L104").Value = 24283.2861
$ws.Range("N104").Value = -29525.2861

$ws.Range("H117").Value = 222225570
$ws.Range("J117").Value = 222225570
$ws.Range("L117").Value = 666676710
$ws.Range("N117").Value = -666683594

$ws.Range("H119").Value = 76936536
$ws.Range("I119").Value = 111119560
$ws.Range("K119").Value = 333358680
$ws.Range("M119").Value = -333353842

$ws.Range("H120").Value = 17307.375
$ws.Range("I120").Value = 6666.3335
$ws.Range("J120").Value = 23692
$ws.Range("K120").Value = 19999.0005
$ws.Range("L120").Value = 71076
$ws.Range("M120").Value = -15161.0005
$ws.Range("N120").Value = -80752

$ws.Range("H121").Value = 2196.9167
$ws.Range("J121").Value = 9999.5
$ws.Range("L121").Value = 29998.5
$ws.Range("N121").Value = -32618.5

$ws.Range("H122").Value = 9525206
$ws.Range("J122").Value = 2331
$ws.Range("L122").Value = 20979
$ws.Range("N122").Value = -25879

$ws.Range("H123").Value = 4248.75
$ws.Range("I123").Value = 4248.75
$ws.Range("K123").Value = 12746.25
$ws.Range("M123").Value = -10296.25

$ws.Range("H124").Value = 3383
$ws.Range("I124").Value = 2872.5
$ws.Range("J124").Value = 3461.5386
$ws.Range("K124").Value = 8617.5
$ws.Range("L124").Value = 10384.6158
$ws.Range("M124").Value = -3707.5
$ws.Range("N124").Value = -20204.6158

$ws.Range("H129").Value = 13082.3
$ws.Range("I129").Value = 1353.6
$ws.Range("J129").Value = 24811
$ws.Range("K129").Value = 4060.8
$ws.Range("L129").Value = 74433
$ws.Range("M129").Value = 939.2000000000003
$ws.Range("N129").Value = -84433

$ws.Range("H132").Value = 1481.0769
$ws.Range("I132").Value = 1028.2222
$ws.Range("K132").Value = 9253.9998
$ws.Range("M132").Value = -6723.9998

$ws.Range("H133").Value = 5115.4
$ws.Range("I133").Value = 5115.4
$ws.Range("K133").Value = 15346.2
$ws.Range("M133").Value = -10286.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 13133.667
$ws.Range("J52").Value = 13133.667
$ws.Range("L52").Value = 13133.667
$ws.Range("N52").Value = -13651.667

$ws.Range("H123").Value = 52666.332
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents() | Out-Null

$ws.Range("H132").Value = 828816.5600000001
$ws.Range("I132").Value = 1171.2222
$ws.Range("K132").Value = 3513.6666
$ws.Range("M132").Value = -983.6665999999996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5988
$ws.Range("J2").Value = 5988
$ws.Range("L2").Value = 5988
$ws.Range("N2").Value = -6212

$ws.Range("H4").Value = 14943.5
$ws.Range("I4").Value = 14943.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 14943.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -14830.5
$ws.Range("N4").ClearContents() | Out-Null

$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 15000
$ws.Range("M25").ClearContents() | Out-Null
$ws.Range("N25").Value = -15460

$ws.Range("H28").Value = 14943.5
$ws.Range("I28").Value = 14943.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 14943.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -14711.5
$ws.Range("N28").ClearContents() | Out-Null

$ws.Range("H37").Value = 14943.5
$ws.Range("I37").Value = 14943.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 14943.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -14836.5
$ws.Range("N37").ClearContents() | Out-Null

$ws.Range("H103").Value = 24444
$ws.Range("J103").Value = 24444
$ws.Range("L103").Value = 24444
$ws.Range("N103").Value = -26788

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents() | Out-Null

$ws.Range("H122").Value = 3465.2083
$ws.Range("I122").Value = 2942.1052
$ws.Range("K122").Value = 8826.3156
$ws.Range("M122").Value = -6376.3156

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 25001882
$ws.Range("I132").Value = 41667840
$ws.Range("J132").Value = 2943.1875
$ws.Range("K132").Value = 125003520
$ws.Range("L132").Value = 8829.5625
$ws.Range("M132").Value = -125000990
$ws.Range("N132").Value = -13889.5625

$ws.Range("H136").Value = 1197.7273
$ws.Range("I136").Value = 750.3333
$ws.Range("J136").Value = 1365.5
$ws.Range("K136").Value = 2250.9999
$ws.Range("L136").Value = 4096.5
$ws.Range("M136").Value = 299.0001000000002
$ws.Range("N136").Value = -9196.5
